$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2 through 497 all currently hold the
# Excel date serial 45171 (2023-09-02). Update them all to 45172 (2023-09-03).
$lastRow = 497
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
